# Added questions total time.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix typo in D19 answer key: 20 -> 10
$ws.Range("D19").Value = "1|2|3|4|5|6|7|8|9|10"

# Add new "Duration (sec)" column
$ws.Range("E1").Value = "Duration (sec)"
$ws.Range("E2").Value = 90
$ws.Range("E3").Value = 90
$ws.Range("E4").Value = 120
$ws.Range("E5").Value = 60
$ws.Range("E6").Value = 60
$ws.Range("E7").Value = 90
$ws.Range("E8").Value = 105
$ws.Range("E9").Value = 90
$ws.Range("E10").Value = 105
$ws.Range("E11").Value = 90
$ws.Range("E12").Value = 120
$ws.Range("E13").Value = 120
$ws.Range("E14").Value = 60
$ws.Range("E15").Value = 60
$ws.Range("E16").Value = 90
$ws.Range("E17").Value = 90
$ws.Range("E18").Value = 90
$ws.Range("E19").Value = 90

# Match header style (bold, centered, wrap) for E1
$ws.Range("E1").Style = $ws.Range("D1").Style

# Match data style (centered, wrap) for E2:E19
$ws.Range("E2:E19").Style = $ws.Range("D2").Style

# Column widths adjustments
$ws.Columns("C").ColumnWidth = 73.54296875
$ws.Columns("D").ColumnWidth = 30.90625
$ws.Columns("E").ColumnWidth = 14.36328125

# Sheet view adjustments
$ws.Range("E11").Select
$excel.ActiveWindow.ScrollColumn = 2
